$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.551.23"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").Value = "4.012.18"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'528.64"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").Value = "'148.73"
$ws.Range("E6").Value = "  +1.92%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "'0.741"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "'44.80"
$ws.Range("E12").Value = "  +3.58%  "
$ws.Range("D13").Value = "'10.80"
$ws.Range("E13").Value = "  +3.51%  "
$ws.Range("D14").Value = "4.650.29"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "4.000.32"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").Value = "'21.39"
$ws.Range("E16").Value = "  +7.60%  "
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  -1.70%  "
$ws.Range("D20").Value = "71.529.42"
$ws.Range("E20").Value = "  +2.88%  "
$ws.Range("D21").Value = "'443.33"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("E22").Value = "  +5.10%  "
$ws.Range("D23").Value = "'93.88"
$ws.Range("E23").Value = "  +6.18%  "
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").Value = "'12.37"
$ws.Range("E25").Value = "  +4.21%  "
$ws.Range("D26").Value = "'4.11"
$ws.Range("E26").Value = "  +5.03%  "
$ws.Range("D27").Value = "'11.10"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'37.22"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").Value = "'705.95"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").Value = "'13.69"
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("D32").Value = "'7.07"
$ws.Range("E32").Value = "  +18.78%  "
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").Value = "'68.37"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "0.0₃0904"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").Value = "'41.14"
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("D38").Value = "'3.63"
$ws.Range("E38").Value = "  +19.75%  "
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'0.0496"
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'2.88"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").Value = "'3.15"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").Value = "'3.56"
$ws.Range("E45").Value = "  +5.50%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.147"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'3.23"
$ws.Range("E47").Value = "  +8.68%  "
$ws.Range("D48").Value = "'0.000285"
$ws.Range("E48").Value = "  +22.14%  "
$ws.Range("E49").Value = "  +6.26%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("E51").Value = "  -3.22%  "
